$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 41334
$ws.Range("H2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("H3").Value = 43070
$ws.Range("H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("H4").Value = 43405
$ws.Range("H4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("H5").Value = 41110
$ws.Range("H5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
